# Week 5 -> Week 6 weekly status report update

$wb = $excel.ActiveWorkbook

# --- Report sheet ---
$report = $wb.Worksheets.Item("Report")

# Status date: 2021-11-23 -> 2021-11-30
$report.Range("C5").Value = 20211130

# "For Period" label: Week 5 -> Week 6
$report.Range("C7").Value = "Week 6"

# Mark task progress for week 6 (columns I/J/K correspond to tasks 6/7/8
# finishing in weeks 6/7/8 respectively; only task 6 (row 27, col I) actually
# completes this week per the source data)
$report.Range("I27").Value = 7500
$report.Range("J28").Value = 7500
$report.Range("K29").Value = 7500

# Extend the weekly totals formula (originally only filled through G35) across
# the newly used H:K columns so the Totals row keeps summing each week.
$report.Range("H35:K35").FormulaR1C1 = "=SUM(R[-13]C:R[-1]C)"

# --- EV (Earned Value) sheet ---
$ev = $wb.Worksheets.Item("EV")

# Tasks 1-5 are 100% complete by week 6; task 6 is 75% complete
$ev.Range("I9").Value = 1
$ev.Range("I10").Value = 1
$ev.Range("I11").Value = 1
$ev.Range("I12").Value = 1
$ev.Range("I13").Value = 1
$ev.Range("I14").Value = 0.75

# --- AC (Actual Cost) sheet ---
$ac = $wb.Worksheets.Item("AC")

# Actual cost incurred during week 6 for task 6
$ac.Range("I14").Value = 7500

# Recalculate the whole workbook so that all dependent formulas (totals,
# chart caches, CPI/SPI/EAC metrics, etc.) are refreshed.
$wb.Application.CalculateFull()
